$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 463554
$ws.Range("I6").Value = 463554
$ws.Range("K6").Value = 1390662
$ws.Range("M6").Value = -1390550

$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 12
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 127
$ws.Range("N8").ClearContents()

$ws.Range("H28").Value = 969.5
$ws.Range("I28").Value = 969.5
$ws.Range("K28").Value = 969.5
$ws.Range("M28").Value = -484.5

$ws.Range("H31").Value = 543.2
$ws.Range("I31").Value = 543.2
$ws.Range("K31").Value = 1629.6
$ws.Range("M31").Value = -1399.6

$ws.Range("H88").Value = 5900
$ws.Range("I88").Value = 3400
$ws.Range("K88").Value = 3400
$ws.Range("M88").Value = -2994

$ws.Range("H91").Value = 5900
$ws.Range("I91").Value = 3400
$ws.Range("K91").Value = 3400
$ws.Range("M91").Value = -1996

$ws.Range("H101").Value = 980
$ws.Range("I101").Value = 980
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2940
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = -1318
$ws.Range("N101").ClearContents()

$ws.Range("H116").Value = 5749
$ws.Range("J116").Value = 5500
$ws.Range("L116").Value = 5500
$ws.Range("N116").Value = -12384

$ws.Range("H118").Value = 220
$ws.Range("I118").Value = 220
$ws.Range("K118").Value = 660
$ws.Range("M118").Value = 997

$ws.Range("H129").Value = 3466.6667
$ws.Range("I129").Value = 1200
$ws.Range("J129").Value = 8000
$ws.Range("K129").Value = 3600
$ws.Range("L129").Value = 24000
$ws.Range("M129").Value = 1400
$ws.Range("N129").Value = -34000

$ws.Range("H132").Value = 5648.727
$ws.Range("I132").Value = 5648.727
$ws.Range("K132").Value = 16946.181
$ws.Range("M132").Value = -14416.181

$ws.Range("H137").Value = 6820.048
$ws.Range("I137").Value = 5689.294
$ws.Range("K137").Value = 17067.882
$ws.Range("M137").Value = -14517.882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 32325.25
$ws.Range("J92").Value = 32325.25
$ws.Range("L92").Value = 32325.25
$ws.Range("N92").Value = -37317.25

$ws.Range("H97").Value = 2432.7778
$ws.Range("I97").Value = 2360.5
$ws.Range("K97").Value = 2360.5
$ws.Range("M97").Value = -1864.5

$ws.Range("H122").Value = 2124.4
$ws.Range("I122").Value = 1905.5
$ws.Range("K122").Value = 5716.5
$ws.Range("M122").Value = -3266.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 475
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 450
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 450
$ws.Range("M7").Value = -387
$ws.Range("N7").Value = -676

$ws.Range("H86").Value = 1541.2
$ws.Range("I86").Value = 1541.2
$ws.Range("K86").Value = 1541.2
$ws.Range("M86").Value = -418.2

$ws.Range("H89").Value = 1541.2
$ws.Range("I89").Value = 1541.2
$ws.Range("K89").Value = 7706
$ws.Range("M89").Value = -2090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5236.6665
$ws.Range("I62").Value = 4105
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 4105
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -3481
$ws.Range("N62").Value = -8748

$ws.Range("H65").Value = 5236.6665
$ws.Range("I65").Value = 4105
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 20525
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -17405
$ws.Range("N65").Value = -43740

$ws.Range("H82").Value = 72498.75
$ws.Range("J82").Value = 72498.75
$ws.Range("L82").Value = 72498.75
$ws.Range("N82").Value = -73220.75

$ws.Range("H85").Value = 72498.75
$ws.Range("J85").Value = 72498.75
$ws.Range("L85").Value = 72498.75
$ws.Range("N85").Value = -74994.75

$ws.Range("H86").Value = 7008
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 7008
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 7008
$ws.Range("M86").Value = -9254

$ws.Range("H89").Value = 7008
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 7008
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 35040
$ws.Range("N89").Value = -46272
$ws.Range("M89").ClearContents()

$ws.Range("H99").Value = 5499

$ws.Range("H107").Value = 6817.8184
$ws.Range("I107").Value = 6817.8184
$ws.Range("K107").Value = 6817.8184
$ws.Range("M107").Value = -4897.8184

$ws.Range("H122").Value = 1598.2222
$ws.Range("I122").Value = 1673.25
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 5019.75
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = -2569.75
$ws.Range("N122").Value = -7894

$ws.Range("H126").Value = 5499

$ws.Range("H132").Value = 3126.2354
$ws.Range("I132").Value = 1626.8462
$ws.Range("K132").Value = 4880.5386
$ws.Range("M132").Value = -2350.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1694.8572
$ws.Range("I6").Value = 1694.8572
$ws.Range("K6").Value = 5084.571599999999
$ws.Range("M6").Value = -4971.571599999999

$ws.Range("H12").Value = 32.933334
$ws.Range("I12").Value = 35.125
$ws.Range("J12").Value = 30.428572
$ws.Range("K12").Value = 105.375
$ws.Range("L12").Value = 91.28571599999999
$ws.Range("M12").Value = 67.625
$ws.Range("N12").Value = -437.285716

$ws.Range("H14").Value = 1459.25
$ws.Range("I14").Value = 1459.25
$ws.Range("K14").Value = 4377.75
$ws.Range("M14").Value = -4204.75

$ws.Range("H52").Value = 1999.5
$ws.Range("J52").Value = 1999.5
$ws.Range("L52").Value = 5998.5
$ws.Range("N52").Value = -6530.5

$ws.Range("H92").Value = 629.6667
$ws.Range("I92").Value = 550
$ws.Range("K92").Value = 1650
$ws.Range("M92").Value = -402

$ws.Range("H98").Value = 877.5
$ws.Range("I98").Value = 1555
$ws.Range("J98").Value = 200
$ws.Range("K98").Value = 4665
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = -3167
$ws.Range("N98").Value = -3596

$ws.Range("H114").Value = 1525.5
$ws.Range("I114").Value = 1145.75
$ws.Range("K114").Value = 3437.25
$ws.Range("M114").Value = -183.25

$ws.Range("H117").Value = 4199.8
$ws.Range("J117").Value = 6000
$ws.Range("L117").Value = 18000
$ws.Range("N117").Value = -24884

$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

$ws.Range("H125").Value = 2600
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 6000
$ws.Range("M125").Value = -1080

$ws.Range("H129").Value = 379
$ws.Range("I129").Value = 455.66666
$ws.Range("J129").Value = 149
$ws.Range("K129").Value = 1366.99998
$ws.Range("L129").Value = 447
$ws.Range("M129").Value = 3633.00002
$ws.Range("N129").Value = -10447

$ws.Range("H140").Value = 1649.75
$ws.Range("I140").Value = 1649.75
$ws.Range("K140").Value = 4949.25
$ws.Range("M140").Value = 230.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 19.461538
$ws.Range("I2").Value = 18.428572
$ws.Range("J2").Value = 20.666666
$ws.Range("K2").Value = 18.428572
$ws.Range("L2").Value = 20.666666
$ws.Range("M2").Value = 94.571428
$ws.Range("N2").Value = -246.666666

$ws.Range("H122").Value = 1797.75
$ws.Range("I122").Value = 1797.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5393.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2943.25
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 10097.6
$ws.Range("J126").Value = 8000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -28940

$ws.Range("H132").Value = 7287.0557
$ws.Range("I132").Value = 5631
$ws.Range("K132").Value = 16893
$ws.Range("M132").Value = -14363

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 462
$ws.Range("I2").Value = 462
$ws.Range("K2").Value = 462
$ws.Range("M2").Value = -350

$ws.Range("H99").Value = 50000
$ws.Range("I99").Value = 50000
$ws.Range("K99").Value = 50000
$ws.Range("M99").Value = -47005

$ws.Range("H132").Value = 7905
$ws.Range("J132").Value = 17000
$ws.Range("L132").Value = 51000
$ws.Range("N132").Value = -56060

$ws.Range("H136").Value = 8944.666999999999
$ws.Range("I136").Value = 3698.25
$ws.Range("K136").Value = 11094.75
$ws.Range("M136").Value = -8544.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 983.5
$ws.Range("I107").Value = 983.5
$ws.Range("K107").Value = 2950.5
$ws.Range("M107").Value = -1030.5

$ws.Range("H136").Value = 8309.333000000001
$ws.Range("I136").Value = 6445.4287
$ws.Range("K136").Value = 19336.2861
$ws.Range("M136").Value = -16786.2861
